$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "3-nov" column (CE) by copying formats from the previous
# date column (CD) and then filling in the values for each row.

# Row 1 (header) - copy format from CD1, then set the new shared string value
$ws.Range("CD1").Copy()
$ws.Range("CE1").PasteSpecial(-4122)
$ws.Range("CE1").Value = "3-nov"

# Row 2
$ws.Range("CD2").Copy()
$ws.Range("CE2").PasteSpecial(-4122)
$ws.Range("CE2").Value = 6

# Row 3
$ws.Range("CD3").Copy()
$ws.Range("CE3").PasteSpecial(-4122)
$ws.Range("CE3").Value = 10

# Row 4
$ws.Range("CD4").Copy()
$ws.Range("CE4").PasteSpecial(-4122)
$ws.Range("CE4").Value = 8

# Row 5
$ws.Range("CD5").Copy()
$ws.Range("CE5").PasteSpecial(-4122)
$ws.Range("CE5").Value = 7

# Row 6 (also gets an extra empty, styled cell in CF6)
$ws.Range("CD6").Copy()
$ws.Range("CE6").PasteSpecial(-4122)
$ws.Range("CE6").Value = 9
$ws.Range("CD6").Copy()
$ws.Range("CF6").PasteSpecial(-4122)

# Row 7
$ws.Range("CD7").Copy()
$ws.Range("CE7").PasteSpecial(-4122)
$ws.Range("CE7").Value = 8

# Row 8
$ws.Range("CD8").Copy()
$ws.Range("CE8").PasteSpecial(-4122)
$ws.Range("CE8").Value = 11

# Row 9
$ws.Range("CD9").Copy()
$ws.Range("CE9").PasteSpecial(-4122)
$ws.Range("CE9").Value = 12

# Row 10
$ws.Range("CD10").Copy()
$ws.Range("CE10").PasteSpecial(-4122)
$ws.Range("CE10").Value = 7

# Row 11
$ws.Range("CD11").Copy()
$ws.Range("CE11").PasteSpecial(-4122)
$ws.Range("CE11").Value = 0

# Update the active selection to match the new last-used cell
$ws.Range("CE1").Select()
